# Add a new "player_fix" worksheet at the end of the workbook containing a
# small lookup table of player name fixes, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# Create the new sheet and name it.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "player_fix"

# Move it to be the last tab in the workbook.
$newSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-acquire the sheet by name (the reference used for Move should not be
# reused for further writes) and populate it with the player name fixes.
$ws = $wb.Worksheets.Item("player_fix")

$ws.Range("A1").Value = "player_name"
$ws.Range("B1").Value = "new_name"
$ws.Range("A2").Value = "Nicolas Claxton"
$ws.Range("B2").Value = "Nic Claxton"
$ws.Range("A3").Value = "OG Anunoby"
$ws.Range("B3").Value = "O.G. Anunoby"

# Make the new sheet the active/selected tab, as in the original workbook.
$ws.Activate()
